$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, and whether the text must be
# force-typed as Text (so Excel does not auto-coerce numeric-looking
# strings like "1.000" or "0.9995" into real numbers).
$updates = @(
    @{ Cell = "D2"; Value = "29.050.12"; ForceText = $false },
    @{ Cell = "E2"; Value = "  +0.03%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "1.830.19"; ForceText = $false },
    @{ Cell = "E3"; Value = "  +0.00%  "; ForceText = $false },
    @{ Cell = "D4"; Value = "0.9988"; ForceText = $true },
    @{ Cell = "E4"; Value = "  +0.01%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "244.84"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +1.52%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "0.6330"; ForceText = $true },
    @{ Cell = "E6"; Value = "  +1.01%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  +0.01%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "0.07545"; ForceText = $true },
    @{ Cell = "E8"; Value = "  -1.04%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "0.2944"; ForceText = $true },
    @{ Cell = "E9"; Value = "  +0.99%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "23.15"; ForceText = $true },
    @{ Cell = "E10"; Value = "  +1.69%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "0.07709"; ForceText = $true },
    @{ Cell = "E11"; Value = "  +0.83%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "1.830.53"; ForceText = $false },
    @{ Cell = "E12"; Value = "  +0.10%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "5.003"; ForceText = $true },
    @{ Cell = "E13"; Value = "  +0.84%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "0.6709"; ForceText = $true },
    @{ Cell = "E14"; Value = "  +0.81%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "83.20"; ForceText = $true },
    @{ Cell = "E15"; Value = "  +0.97%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "0.000009530"; ForceText = $true },
    @{ Cell = "E16"; Value = "  +1.99%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "6.084"; ForceText = $true },
    @{ Cell = "E17"; Value = "  +1.52%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "29.079.38"; ForceText = $false },
    @{ Cell = "E18"; Value = "  +0.75%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "12.59"; ForceText = $true },
    @{ Cell = "E19"; Value = "  +2.13%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "227.03"; ForceText = $true },
    @{ Cell = "E20"; Value = "  +0.94%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "0.9993"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -0.06%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "7.155"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -0.74%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "1.000"; ForceText = $true },
    @{ Cell = "E23"; Value = "  -0.01%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "160.11"; ForceText = $true },
    @{ Cell = "E24"; Value = "  +0.23%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "0.1431"; ForceText = $true },
    @{ Cell = "E25"; Value = "  +5.08%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "8.517"; ForceText = $true },
    @{ Cell = "E26"; Value = "  +1.13%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "17.95"; ForceText = $true },
    @{ Cell = "E27"; Value = "  +0.76%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "1.506"; ForceText = $true },
    @{ Cell = "E28"; Value = "  +0.90%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "4.148"; ForceText = $true },
    @{ Cell = "E29"; Value = "  +2.37%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "4.078"; ForceText = $true },
    @{ Cell = "E30"; Value = "  +1.22%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "0.05490"; ForceText = $true },
    @{ Cell = "E31"; Value = "  +5.52%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "1.201"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -0.27%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "1.860"; ForceText = $true },
    @{ Cell = "E33"; Value = "  +0.63%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "0.7453"; ForceText = $true },
    @{ Cell = "E34"; Value = "  +1.88%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "1.141"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -1.14%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "2.658"; ForceText = $true },
    @{ Cell = "E36"; Value = "  +1.74%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "1.244.96"; ForceText = $false },
    @{ Cell = "E37"; Value = "  -2.48%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "2.759"; ForceText = $true },
    @{ Cell = "E38"; Value = "  +0.00%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "0.01787"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -0.02%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "6.605"; ForceText = $true },
    @{ Cell = "E40"; Value = "  +1.22%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "0.9022"; ForceText = $true },
    @{ Cell = "E41"; Value = "  +1.35%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "0.9995"; ForceText = $true },
    @{ Cell = "E42"; Value = "  -0.02%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "101.44"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -0.06%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "1.982.05"; ForceText = $false },
    @{ Cell = "E44"; Value = "  +0.35%  "; ForceText = $false },
    @{ Cell = "B45"; Value = "Aave"; ForceText = $false },
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; ForceText = $false },
    @{ Cell = "D45"; Value = "65.09"; ForceText = $true },
    @{ Cell = "E45"; Value = "  +2.10%  "; ForceText = $false },
    @{ Cell = "B46"; Value = "BabyDogeCoin"; ForceText = $false },
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; ForceText = $false },
    @{ Cell = "D46"; Value = "0.00000000123"; ForceText = $true },
    @{ Cell = "E46"; Value = "  +1.91%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "0.5099"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -0.15%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "0.4070"; ForceText = $true },
    @{ Cell = "E48"; Value = "  +2.32%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "9.001"; ForceText = $true },
    @{ Cell = "E49"; Value = "  +1.82%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "1.655"; ForceText = $true },
    @{ Cell = "E50"; Value = "  +0.73%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "6.784"; ForceText = $true },
    @{ Cell = "E51"; Value = "  +1.29%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Mark as Text so Excel keeps the exact string instead of parsing it as a
        # number (e.g. "1.000" -> 1, "0.00000000123" -> 1.23E-9).
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        # Drop the temporary Text format again so the cell keeps its original
        # (default) style, matching the source workbook exactly.
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}

Write-Host "Applied cryptos update"